# Rename 'variable' and 'long_name' to 'variable-code' and 'variable-label'
# Close #144
#
# - "Variables" sheet header row: variable -> variable-code,
#   da_long_name/en_long_name/kl_long_name -> da_variable-label/en_variable-label/kl_variable-label
# - "Codelists" sheet header row: variable -> variable-code,
#   da_code_label/en_code_label/kl_code_label -> da_code-label/en_code-label/kl_code-label
# - Active sheet / selection changes: Codelists becomes the active/selected tab
#   (A2 selected), Variables is no longer the selected tab (D1 selected there).

$wb = $excel.ActiveWorkbook

$wsVariables = $wb.Worksheets.Item("Variables")
$wsCodelists = $wb.Worksheets.Item("Codelists")

# --- Rename header labels on the "Variables" sheet ---
$wsVariables.Range("C1").Value = "variable-code"
$wsVariables.Range("D1").Value = "da_variable-label"
$wsVariables.Range("E1").Value = "en_variable-label"
$wsVariables.Range("F1").Value = "kl_variable-label"

# --- Rename header labels on the "Codelists" sheet ---
$wsCodelists.Range("A1").Value = "variable-code"
$wsCodelists.Range("D1").Value = "da_code-label"
$wsCodelists.Range("E1").Value = "en_code-label"
$wsCodelists.Range("F1").Value = "kl_code-label"

# --- Update view/selection state to match final workbook state ---
# Variables sheet: no longer the active tab, selection moves to D1
$wsVariables.Activate()
$wsVariables.Range("D1").Select()

# Codelists sheet: becomes the active tab, selection moves to A2
$wsCodelists.Activate()
$wsCodelists.Range("A2").Select()
